# Auto-generated edit script
# Applies the diff: updates weekly price records for rows 692-756 (shift by one
# week / two rows) and appends two brand-new rows (757, 758) for the latest week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 692-756 with the shifted weekly data ---
$ws.Cells.Item(692, 4).Value = 44769
$ws.Cells.Item(692, 10).Value = 2500
$ws.Cells.Item(692, 11).Value = 850
$ws.Cells.Item(692, 12).Value = 900
$ws.Cells.Item(692, 13).Value = 880
$ws.Cells.Item(692, 16).Value = 220
$ws.Cells.Item(693, 4).Value = 44769
$ws.Cells.Item(693, 10).Value = 1100
$ws.Cells.Item(693, 11).Value = 600
$ws.Cells.Item(693, 12).Value = 600
$ws.Cells.Item(693, 13).Value = 600
$ws.Cells.Item(693, 16).Value = 150
$ws.Cells.Item(694, 9).Value = "Primera"
$ws.Cells.Item(694, 10).Value = 1200
$ws.Cells.Item(694, 11).Value = 550
$ws.Cells.Item(694, 12).Value = 550
$ws.Cells.Item(694, 13).Value = 550
$ws.Cells.Item(694, 16).Value = 138
$ws.Cells.Item(695, 4).Value = 44596
$ws.Cells.Item(695, 9).Value = "Segunda"
$ws.Cells.Item(695, 10).Value = 950
$ws.Cells.Item(695, 11).Value = 400
$ws.Cells.Item(695, 12).Value = 400
$ws.Cells.Item(695, 13).Value = 400
$ws.Cells.Item(695, 16).Value = 100
$ws.Cells.Item(696, 4).Value = 44596
$ws.Cells.Item(696, 9).Value = "Tercera"
$ws.Cells.Item(696, 10).Value = 850
$ws.Cells.Item(697, 4).Value = 44340
$ws.Cells.Item(697, 10).Value = 1100
$ws.Cells.Item(697, 11).Value = 600
$ws.Cells.Item(697, 12).Value = 600
$ws.Cells.Item(697, 13).Value = 600
$ws.Cells.Item(697, 16).Value = 150
$ws.Cells.Item(698, 4).Value = 44340
$ws.Cells.Item(698, 10).Value = 1000
$ws.Cells.Item(698, 11).Value = 500
$ws.Cells.Item(698, 12).Value = 500
$ws.Cells.Item(698, 13).Value = 500
$ws.Cells.Item(698, 16).Value = 125
$ws.Cells.Item(699, 4).Value = 44496
$ws.Cells.Item(699, 10).Value = 2400
$ws.Cells.Item(699, 11).Value = 500
$ws.Cells.Item(699, 13).Value = 525
$ws.Cells.Item(699, 16).Value = 131
$ws.Cells.Item(700, 4).Value = 44496
$ws.Cells.Item(700, 11).Value = 400
$ws.Cells.Item(700, 12).Value = 400
$ws.Cells.Item(700, 13).Value = 400
$ws.Cells.Item(700, 16).Value = 100
$ws.Cells.Item(701, 4).Value = 44399
$ws.Cells.Item(701, 10).Value = 1600
$ws.Cells.Item(701, 12).Value = 550
$ws.Cells.Item(701, 13).Value = 550
$ws.Cells.Item(701, 16).Value = 138
$ws.Cells.Item(702, 4).Value = 44399
$ws.Cells.Item(702, 10).Value = 1300
$ws.Cells.Item(702, 11).Value = 450
$ws.Cells.Item(702, 12).Value = 450
$ws.Cells.Item(702, 13).Value = 450
$ws.Cells.Item(702, 16).Value = 112
$ws.Cells.Item(703, 4).Value = 44425
$ws.Cells.Item(703, 10).Value = 3300
$ws.Cells.Item(703, 11).Value = 550
$ws.Cells.Item(703, 12).Value = 600
$ws.Cells.Item(703, 13).Value = 577
$ws.Cells.Item(703, 16).Value = 144
$ws.Cells.Item(704, 4).Value = 44425
$ws.Cells.Item(704, 10).Value = 1600
$ws.Cells.Item(705, 4).Value = 44377
$ws.Cells.Item(705, 10).Value = 2100
$ws.Cells.Item(705, 12).Value = 500
$ws.Cells.Item(705, 13).Value = 500
$ws.Cells.Item(705, 16).Value = 125
$ws.Cells.Item(706, 4).Value = 44377
$ws.Cells.Item(706, 10).Value = 2600
$ws.Cells.Item(707, 4).Value = 44512
$ws.Cells.Item(707, 10).Value = 2900
$ws.Cells.Item(707, 12).Value = 550
$ws.Cells.Item(707, 13).Value = 524
$ws.Cells.Item(707, 16).Value = 131
$ws.Cells.Item(708, 4).Value = 44512
$ws.Cells.Item(708, 10).Value = 1600
$ws.Cells.Item(709, 4).Value = 44397
$ws.Cells.Item(709, 10).Value = 1300
$ws.Cells.Item(709, 11).Value = 500
$ws.Cells.Item(709, 12).Value = 500
$ws.Cells.Item(709, 13).Value = 500
$ws.Cells.Item(709, 16).Value = 125
$ws.Cells.Item(710, 4).Value = 44397
$ws.Cells.Item(710, 11).Value = 400
$ws.Cells.Item(710, 12).Value = 400
$ws.Cells.Item(710, 13).Value = 400
$ws.Cells.Item(710, 16).Value = 100
$ws.Cells.Item(711, 4).Value = 44181
$ws.Cells.Item(711, 10).Value = 2050
$ws.Cells.Item(711, 11).Value = 400
$ws.Cells.Item(711, 12).Value = 450
$ws.Cells.Item(711, 13).Value = 429
$ws.Cells.Item(711, 16).Value = 107
$ws.Cells.Item(712, 4).Value = 44181
$ws.Cells.Item(712, 10).Value = 1200
$ws.Cells.Item(712, 11).Value = 350
$ws.Cells.Item(712, 12).Value = 350
$ws.Cells.Item(712, 13).Value = 350
$ws.Cells.Item(712, 16).Value = 88
$ws.Cells.Item(713, 4).Value = 44497
$ws.Cells.Item(713, 10).Value = 2650
$ws.Cells.Item(713, 11).Value = 500
$ws.Cells.Item(713, 12).Value = 550
$ws.Cells.Item(713, 13).Value = 525
$ws.Cells.Item(713, 16).Value = 131
$ws.Cells.Item(714, 4).Value = 44497
$ws.Cells.Item(714, 10).Value = 1280
$ws.Cells.Item(714, 11).Value = 400
$ws.Cells.Item(714, 12).Value = 400
$ws.Cells.Item(714, 13).Value = 400
$ws.Cells.Item(714, 16).Value = 100
$ws.Cells.Item(715, 4).Value = 44285
$ws.Cells.Item(715, 10).Value = 1800
$ws.Cells.Item(715, 11).Value = 700
$ws.Cells.Item(715, 12).Value = 700
$ws.Cells.Item(715, 13).Value = 700
$ws.Cells.Item(715, 16).Value = 175
$ws.Cells.Item(716, 4).Value = 44285
$ws.Cells.Item(716, 10).Value = 1200
$ws.Cells.Item(716, 11).Value = 600
$ws.Cells.Item(716, 12).Value = 600
$ws.Cells.Item(716, 13).Value = 600
$ws.Cells.Item(716, 16).Value = 150
$ws.Cells.Item(717, 4).Value = 44362
$ws.Cells.Item(717, 10).Value = 3200
$ws.Cells.Item(717, 11).Value = 550
$ws.Cells.Item(717, 13).Value = 575
$ws.Cells.Item(717, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(717, 16).Value = 144
$ws.Cells.Item(717, 17).Value = 4
$ws.Cells.Item(718, 4).Value = 44362
$ws.Cells.Item(718, 10).Value = 1800
$ws.Cells.Item(718, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(718, 16).Value = 100
$ws.Cells.Item(718, 17).Value = 4
$ws.Cells.Item(719, 4).Value = 44557
$ws.Cells.Item(719, 10).Value = 2250
$ws.Cells.Item(719, 11).Value = 500
$ws.Cells.Item(719, 12).Value = 600
$ws.Cells.Item(719, 13).Value = 544
$ws.Cells.Item(719, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(719, 16).Value = 45
$ws.Cells.Item(719, 17).Value = 12
$ws.Cells.Item(720, 4).Value = 44557
$ws.Cells.Item(720, 10).Value = 1700
$ws.Cells.Item(720, 11).Value = 400
$ws.Cells.Item(720, 12).Value = 400
$ws.Cells.Item(720, 13).Value = 400
$ws.Cells.Item(720, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(720, 16).Value = 33
$ws.Cells.Item(720, 17).Value = 12
$ws.Cells.Item(721, 4).Value = 44747
$ws.Cells.Item(721, 10).Value = 2500
$ws.Cells.Item(721, 11).Value = 650
$ws.Cells.Item(721, 12).Value = 700
$ws.Cells.Item(721, 13).Value = 676
$ws.Cells.Item(721, 16).Value = 169
$ws.Cells.Item(722, 4).Value = 44747
$ws.Cells.Item(722, 11).Value = 500
$ws.Cells.Item(722, 12).Value = 500
$ws.Cells.Item(722, 13).Value = 500
$ws.Cells.Item(722, 16).Value = 125
$ws.Cells.Item(723, 4).Value = 44357
$ws.Cells.Item(723, 10).Value = 2900
$ws.Cells.Item(723, 11).Value = 500
$ws.Cells.Item(723, 13).Value = 524
$ws.Cells.Item(723, 16).Value = 131
$ws.Cells.Item(724, 4).Value = 44357
$ws.Cells.Item(724, 10).Value = 1600
$ws.Cells.Item(724, 11).Value = 400
$ws.Cells.Item(724, 12).Value = 400
$ws.Cells.Item(724, 13).Value = 400
$ws.Cells.Item(724, 16).Value = 100
$ws.Cells.Item(725, 4).Value = 44279
$ws.Cells.Item(725, 10).Value = 1600
$ws.Cells.Item(725, 11).Value = 550
$ws.Cells.Item(725, 12).Value = 550
$ws.Cells.Item(725, 13).Value = 550
$ws.Cells.Item(725, 16).Value = 138
$ws.Cells.Item(726, 4).Value = 44279
$ws.Cells.Item(726, 10).Value = 850
$ws.Cells.Item(726, 11).Value = 450
$ws.Cells.Item(726, 12).Value = 450
$ws.Cells.Item(726, 13).Value = 450
$ws.Cells.Item(726, 16).Value = 112
$ws.Cells.Item(727, 4).Value = 44551
$ws.Cells.Item(727, 10).Value = 180
$ws.Cells.Item(727, 12).Value = 500
$ws.Cells.Item(727, 13).Value = 500
$ws.Cells.Item(727, 16).Value = 125
$ws.Cells.Item(728, 4).Value = 44551
$ws.Cells.Item(728, 10).Value = 160
$ws.Cells.Item(729, 4).Value = 44517
$ws.Cells.Item(729, 10).Value = 3000
$ws.Cells.Item(729, 11).Value = 500
$ws.Cells.Item(729, 12).Value = 550
$ws.Cells.Item(729, 13).Value = 525
$ws.Cells.Item(729, 16).Value = 131
$ws.Cells.Item(730, 4).Value = 44517
$ws.Cells.Item(730, 10).Value = 1600
$ws.Cells.Item(730, 11).Value = 400
$ws.Cells.Item(730, 12).Value = 400
$ws.Cells.Item(730, 13).Value = 400
$ws.Cells.Item(730, 16).Value = 100
$ws.Cells.Item(731, 4).Value = 44757
$ws.Cells.Item(731, 10).Value = 3400
$ws.Cells.Item(731, 11).Value = 800
$ws.Cells.Item(731, 12).Value = 900
$ws.Cells.Item(731, 13).Value = 847
$ws.Cells.Item(731, 16).Value = 212
$ws.Cells.Item(732, 4).Value = 44757
$ws.Cells.Item(732, 10).Value = 1700
$ws.Cells.Item(732, 11).Value = 600
$ws.Cells.Item(732, 12).Value = 600
$ws.Cells.Item(732, 13).Value = 600
$ws.Cells.Item(732, 16).Value = 150
$ws.Cells.Item(733, 4).Value = 44547
$ws.Cells.Item(733, 10).Value = 3000
$ws.Cells.Item(733, 11).Value = 450
$ws.Cells.Item(733, 12).Value = 500
$ws.Cells.Item(733, 13).Value = 477
$ws.Cells.Item(733, 16).Value = 119
$ws.Cells.Item(734, 4).Value = 44547
$ws.Cells.Item(734, 10).Value = 2500
$ws.Cells.Item(734, 11).Value = 350
$ws.Cells.Item(734, 12).Value = 400
$ws.Cells.Item(734, 13).Value = 374
$ws.Cells.Item(734, 16).Value = 94
$ws.Cells.Item(735, 4).Value = 44321
$ws.Cells.Item(735, 10).Value = 1300
$ws.Cells.Item(735, 11).Value = 700
$ws.Cells.Item(735, 12).Value = 700
$ws.Cells.Item(735, 13).Value = 700
$ws.Cells.Item(735, 16).Value = 175
$ws.Cells.Item(736, 4).Value = 44321
$ws.Cells.Item(736, 10).Value = 1850
$ws.Cells.Item(736, 11).Value = 500
$ws.Cells.Item(736, 12).Value = 500
$ws.Cells.Item(736, 13).Value = 500
$ws.Cells.Item(736, 16).Value = 125
$ws.Cells.Item(737, 4).Value = 44438
$ws.Cells.Item(737, 11).Value = 600
$ws.Cells.Item(737, 12).Value = 650
$ws.Cells.Item(737, 13).Value = 626
$ws.Cells.Item(737, 16).Value = 156
$ws.Cells.Item(738, 4).Value = 44438
$ws.Cells.Item(738, 10).Value = 1800
$ws.Cells.Item(739, 4).Value = 44355
$ws.Cells.Item(739, 10).Value = 3100
$ws.Cells.Item(739, 11).Value = 500
$ws.Cells.Item(739, 12).Value = 550
$ws.Cells.Item(739, 13).Value = 521
$ws.Cells.Item(739, 16).Value = 130
$ws.Cells.Item(740, 4).Value = 44355
$ws.Cells.Item(740, 10).Value = 1600
$ws.Cells.Item(741, 4).Value = 44657
$ws.Cells.Item(741, 10).Value = 2800
$ws.Cells.Item(741, 11).Value = 550
$ws.Cells.Item(741, 12).Value = 600
$ws.Cells.Item(741, 13).Value = 571
$ws.Cells.Item(741, 16).Value = 143
$ws.Cells.Item(742, 4).Value = 44657
$ws.Cells.Item(742, 10).Value = 1390
$ws.Cells.Item(743, 4).Value = 44391
$ws.Cells.Item(743, 10).Value = 3100
$ws.Cells.Item(744, 4).Value = 44391
$ws.Cells.Item(744, 10).Value = 1400
$ws.Cells.Item(745, 4).Value = 44453
$ws.Cells.Item(745, 10).Value = 2700
$ws.Cells.Item(745, 11).Value = 500
$ws.Cells.Item(745, 12).Value = 550
$ws.Cells.Item(745, 13).Value = 526
$ws.Cells.Item(745, 16).Value = 132
$ws.Cells.Item(746, 4).Value = 44453
$ws.Cells.Item(746, 11).Value = 400
$ws.Cells.Item(746, 12).Value = 400
$ws.Cells.Item(746, 13).Value = 400
$ws.Cells.Item(746, 16).Value = 100
$ws.Cells.Item(747, 4).Value = 44186
$ws.Cells.Item(747, 10).Value = 2600
$ws.Cells.Item(747, 11).Value = 450
$ws.Cells.Item(747, 12).Value = 500
$ws.Cells.Item(747, 13).Value = 475
$ws.Cells.Item(747, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(747, 16).Value = 119
$ws.Cells.Item(748, 4).Value = 44186
$ws.Cells.Item(748, 10).Value = 1200
$ws.Cells.Item(748, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(749, 4).Value = 44189
$ws.Cells.Item(749, 10).Value = 3000
$ws.Cells.Item(749, 11).Value = 500
$ws.Cells.Item(749, 12).Value = 550
$ws.Cells.Item(749, 13).Value = 527
$ws.Cells.Item(749, 15).Value = "Hijuelas"
$ws.Cells.Item(749, 16).Value = 132
$ws.Cells.Item(750, 4).Value = 44189
$ws.Cells.Item(750, 10).Value = 900
$ws.Cells.Item(750, 11).Value = 350
$ws.Cells.Item(750, 12).Value = 350
$ws.Cells.Item(750, 13).Value = 350
$ws.Cells.Item(750, 15).Value = "Hijuelas"
$ws.Cells.Item(750, 16).Value = 88
$ws.Cells.Item(751, 4).Value = 44609
$ws.Cells.Item(751, 10).Value = 2630
$ws.Cells.Item(751, 11).Value = 550
$ws.Cells.Item(751, 12).Value = 600
$ws.Cells.Item(751, 13).Value = 576
$ws.Cells.Item(751, 16).Value = 144
$ws.Cells.Item(752, 4).Value = 44609
$ws.Cells.Item(752, 10).Value = 1250
$ws.Cells.Item(753, 4).Value = 44489
$ws.Cells.Item(753, 10).Value = 3200
$ws.Cells.Item(753, 11).Value = 500
$ws.Cells.Item(753, 12).Value = 550
$ws.Cells.Item(753, 13).Value = 525
$ws.Cells.Item(753, 16).Value = 131
$ws.Cells.Item(754, 4).Value = 44489
$ws.Cells.Item(754, 10).Value = 1800
$ws.Cells.Item(755, 4).Value = 44358
$ws.Cells.Item(755, 10).Value = 3100
$ws.Cells.Item(755, 11).Value = 550
$ws.Cells.Item(755, 12).Value = 600
$ws.Cells.Item(755, 13).Value = 579
$ws.Cells.Item(755, 16).Value = 145
$ws.Cells.Item(756, 4).Value = 44358
$ws.Cells.Item(756, 10).Value = 1200

# --- Append two brand-new rows (757-758) for the newest price entries ---
$ws.Cells.Item(757, 1).Value = 3
$ws.Cells.Item(757, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(757, 3).Value = "Coquimbo"
$ws.Cells.Item(757, 4).Value = 44572
$ws.Cells.Item(757, 5).Value = 5
$ws.Cells.Item(757, 6).Value = 100114014
$ws.Cells.Item(757, 7).Value = "Betarraga"
$ws.Cells.Item(757, 8).Value = "Sin especificar"
$ws.Cells.Item(757, 9).Value = "Primera"
$ws.Cells.Item(757, 10).Value = 4000
$ws.Cells.Item(757, 11).Value = 500
$ws.Cells.Item(757, 12).Value = 550
$ws.Cells.Item(757, 13).Value = 515
$ws.Cells.Item(757, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(757, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(757, 16).Value = 129
$ws.Cells.Item(757, 17).Value = 4
$ws.Cells.Item(757, 18).Value = "Hortaliza"
$ws.Cells.Item(757, 4).NumberFormat = $ws.Cells.Item(756, 4).NumberFormat

$ws.Cells.Item(758, 1).Value = 3
$ws.Cells.Item(758, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(758, 3).Value = "Coquimbo"
$ws.Cells.Item(758, 4).Value = 44572
$ws.Cells.Item(758, 5).Value = 5
$ws.Cells.Item(758, 6).Value = 100114014
$ws.Cells.Item(758, 7).Value = "Betarraga"
$ws.Cells.Item(758, 8).Value = "Sin especificar"
$ws.Cells.Item(758, 9).Value = "Segunda"
$ws.Cells.Item(758, 10).Value = 1300
$ws.Cells.Item(758, 11).Value = 400
$ws.Cells.Item(758, 12).Value = 400
$ws.Cells.Item(758, 13).Value = 400
$ws.Cells.Item(758, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(758, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(758, 16).Value = 100
$ws.Cells.Item(758, 17).Value = 4
$ws.Cells.Item(758, 18).Value = "Hortaliza"
$ws.Cells.Item(758, 4).NumberFormat = $ws.Cells.Item(756, 4).NumberFormat

